$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.626383
$ws.Range("H2").Value = 10.879149
$ws.Range("I2").Value = 0.04794321320295086
$ws.Range("J2").Value = 0.05075107175358264
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 147.0592853333333
$ws.Range("N2").Value = 441.177856
$ws.Range("O2").Value = 0.9129893958419346
$ws.Range("P2").Value = 0.9274576550077637
$ws.Range("Q2").Value = 533.2932923249494
$ws.Range("R2").Value = 4799.639630924544
$ws.Range("S2").Value = 0.04377164525688317
$ws.Range("T2").Value = 0.0470694699977085

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.626383
$ws.Range("H3").Value = 10.879149
$ws.Range("I3").Value = 0.04794321320295086
$ws.Range("J3").Value = 0.05075107175358264
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3688046666666667
$ws.Range("N3").Value = 1.106414
$ws.Range("O3").Value = 0.00228965310854373
$ws.Range("P3").Value = 0.002325937532793486
$ws.Range("Q3").Value = 1.337426973520667
$ws.Range("R3").Value = 12.036842761686
$ws.Range("S3").Value = 0.0001097733271437112
$ws.Range("T3").Value = 0.0001180438226211532

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.626383
$ws.Range("H4").Value = 10.879149
$ws.Range("I4").Value = 0.04794321320295086
$ws.Range("J4").Value = 0.05075107175358264
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.963579
$ws.Range("N4").Value = 14.890737
$ws.Range("O4").Value = 0.03081542917981618
$ws.Range("P4").Value = 0.03130376520837289
$ws.Range("Q4").Value = 17.999838504757
$ws.Range("R4").Value = 161.998546542813
$ws.Range("S4").Value = 0.00147739069110836
$ws.Range("T4").Value = 0.001588699634247436

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.626383
$ws.Range("H5").Value = 10.879149
$ws.Range("I5").Value = 0.04794321320295086
$ws.Range("J5").Value = 0.05075107175358264
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.14456
$ws.Range("N5").Value = 3.43368
$ws.Range("O5").Value = 0.007105781457704291
$ws.Range("P5").Value = 0.007218387680924443
$ws.Range("Q5").Value = 4.15061292648
$ws.Range("R5").Value = 37.35551633832
$ws.Range("S5").Value = 0.0003406739954002918
$ws.Range("T5").Value = 0.0003663409111397734

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.626383
$ws.Range("H6").Value = 10.879149
$ws.Range("I6").Value = 0.04794321320295086
$ws.Range("J6").Value = 0.05075107175358264
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 7.5382435
$ws.Range("N6").Value = 15.076487
$ws.Range("O6").Value = 0.04679974041200103
$ws.Range("P6").Value = 0.0316942545701456
$ws.Range("Q6").Value = 27.3365580782605
$ws.Range("R6").Value = 164.019348469563
$ws.Range("S6").Value = 0.002243729932415321
$ws.Range("T6").Value = 0.001608517387865773

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 26.85202466666667
$ws.Range("H7").Value = 80.556074
$ws.Range("I7").Value = 0.3550017589220155
$ws.Range("J7").Value = 0.3757929128244233
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 147.0592853333333
$ws.Range("N7").Value = 441.177856
$ws.Range("O7").Value = 0.9129893958419346
$ws.Range("P7").Value = 0.9274576550077637
$ws.Range("Q7").Value = 3948.839557233038
$ws.Range("R7").Value = 35539.55601509735
$ws.Range("S7").Value = 0.3241128414010351
$ws.Range("T7").Value = 0.3485320136966766

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 26.85202466666667
$ws.Range("H8").Value = 80.556074
$ws.Range("I8").Value = 0.3550017589220155
$ws.Range("J8").Value = 0.3757929128244233
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.3688046666666667
$ws.Range("N8").Value = 1.106414
$ws.Range("O8").Value = 0.00228965310854373
$ws.Range("P8").Value = 0.002325937532793486
$ws.Range("Q8").Value = 9.90315200651511
$ws.Range("R8").Value = 89.128368058636
$ws.Range("S8").Value = 0.0008128308808542846
$ws.Range("T8").Value = 0.0008740708404961168

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 26.85202466666667
$ws.Range("H9").Value = 80.556074
$ws.Range("I9").Value = 0.3550017589220155
$ws.Range("J9").Value = 0.3757929128244233
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.963579
$ws.Range("N9").Value = 14.890737
$ws.Range("O9").Value = 0.03081542917981618
$ws.Range("P9").Value = 0.03130376520837289
$ws.Range("Q9").Value = 133.2821457429487
$ws.Range("R9").Value = 1199.539311686538
$ws.Range("S9").Value = 0.01093953156077155
$ws.Range("T9").Value = 0.01176373311002629

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 26.85202466666667
$ws.Range("H10").Value = 80.556074
$ws.Range("I10").Value = 0.3550017589220155
$ws.Range("J10").Value = 0.3757929128244233
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.14456
$ws.Range("N10").Value = 3.43368
$ws.Range("O10").Value = 0.007105781457704291
$ws.Range("P10").Value = 0.007218387680924443
$ws.Range("Q10").Value = 30.73375335248
$ws.Range("R10").Value = 276.6037801723199
$ws.Range("S10").Value = 0.002522564916000467
$ws.Range("T10").Value = 0.002712618932510531

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 26.85202466666667
$ws.Range("H11").Value = 80.556074
$ws.Range("I11").Value = 0.3550017589220155
$ws.Range("J11").Value = 0.3757929128244233
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 7.5382435
$ws.Range("N11").Value = 15.076487
$ws.Range("O11").Value = 0.04679974041200103
$ws.Range("P11").Value = 0.0316942545701456
$ws.Range("Q11").Value = 202.4171004053397
$ws.Range("R11").Value = 1214.502602432038
$ws.Range("S11").Value = 0.0166139901633541
$ws.Range("T11").Value = 0.0119104762447138

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 13.189183
$ws.Range("H12").Value = 39.567549
$ws.Range("I12").Value = 0.1743698369812937
$ws.Range("J12").Value = 0.1845820402324113
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 147.0592853333333
$ws.Range("N12").Value = 441.177856
$ws.Range("O12").Value = 0.9129893958419346
$ws.Range("P12").Value = 0.9274576550077637
$ws.Range("Q12").Value = 1939.59182611055
$ws.Range("R12").Value = 17456.32643499495
$ws.Range("S12").Value = 0.159197812118608
$ws.Range("T12").Value = 0.1711920261905008

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 13.189183
$ws.Range("H13").Value = 39.567549
$ws.Range("I13").Value = 0.1743698369812937
$ws.Range("J13").Value = 0.1845820402324113
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.3688046666666667
$ws.Range("N13").Value = 1.106414
$ws.Range("O13").Value = 0.00228965310854373
$ws.Range("P13").Value = 0.002325937532793486
$ws.Range("Q13").Value = 4.864232239920667
$ws.Range("R13").Value = 43.778090159286
$ws.Range("S13").Value = 0.0003992464392804827
$ws.Range("T13").Value = 0.0004293262952561626

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.189183
$ws.Range("H14").Value = 39.567549
$ws.Range("I14").Value = 0.1743698369812937
$ws.Range("J14").Value = 0.1845820402324113
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.963579
$ws.Range("N14").Value = 14.890737
$ws.Range("O14").Value = 0.03081542917981618
$ws.Range("P14").Value = 0.03130376520837289
$ws.Range("Q14").Value = 65.465551765957
$ws.Range("R14").Value = 589.189965893613
$ws.Range("S14").Value = 0.00537328136259315
$ws.Range("T14").Value = 0.00577811284911784

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.189183
$ws.Range("H15").Value = 39.567549
$ws.Range("I15").Value = 0.1743698369812937
$ws.Range("J15").Value = 0.1845820402324113
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.14456
$ws.Range("N15").Value = 3.43368
$ws.Range("O15").Value = 0.007105781457704291
$ws.Range("P15").Value = 0.007218387680924443
$ws.Range("Q15").Value = 15.09581129448
$ws.Range("R15").Value = 135.86230165032
$ws.Range("S15").Value = 0.001239033954404597
$ws.Range("T15").Value = 0.001332384725333537

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.189183
$ws.Range("H16").Value = 39.567549
$ws.Range("I16").Value = 0.1743698369812937
$ws.Range("J16").Value = 0.1845820402324113
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 7.5382435
$ws.Range("N16").Value = 15.076487
$ws.Range("O16").Value = 0.04679974041200103
$ws.Range("P16").Value = 0.0316942545701456
$ws.Range("Q16").Value = 99.4232730200605
$ws.Range("R16").Value = 596.539638120363
$ws.Range("S16").Value = 0.008160463106407485
$ws.Range("T16").Value = 0.005850190172202899

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.41708766666666
$ws.Range("H17").Value = 58.25126299999999
$ws.Range("I17").Value = 0.2567069098281641
$ws.Range("J17").Value = 0.2717412941260215
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 147.0592853333333
$ws.Range("N17").Value = 441.177856
$ws.Range("O17").Value = 0.9129893958419346
$ws.Range("P17").Value = 0.9274576550077637
$ws.Range("Q17").Value = 2855.463035514681
$ws.Range("R17").Value = 25699.16731963213
$ws.Range("S17").Value = 0.2343706865124656
$ws.Range("T17").Value = 0.2520285434188949

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 19.41708766666666
$ws.Range("H18").Value = 58.25126299999999
$ws.Range("I18").Value = 0.2567069098281641
$ws.Range("J18").Value = 0.2717412941260215
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.3688046666666667
$ws.Range("N18").Value = 1.106414
$ws.Range("O18").Value = 0.00228965310854373
$ws.Range("P18").Value = 0.002325937532793486
$ws.Range("Q18").Value = 7.161112544542443
$ws.Range("R18").Value = 64.450012900882
$ws.Range("S18").Value = 0.000587769774072711
$ws.Range("T18").Value = 0.0006320532752175876

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 19.41708766666666
$ws.Range("H19").Value = 58.25126299999999
$ws.Range("I19").Value = 0.2567069098281641
$ws.Range("J19").Value = 0.2717412941260215
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 4.963579
$ws.Range("N19").Value = 14.890737
$ws.Range("O19").Value = 0.03081542917981618
$ws.Range("P19").Value = 0.03130376520837289
$ws.Range("Q19").Value = 96.37824858342566
$ws.Range("R19").Value = 867.4042372508309
$ws.Range("S19").Value = 0.00791053359977925
$ws.Range("T19").Value = 0.008506525668740376

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 19.41708766666666
$ws.Range("H20").Value = 58.25126299999999
$ws.Range("I20").Value = 0.2567069098281641
$ws.Range("J20").Value = 0.2717412941260215
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 1.14456
$ws.Range("N20").Value = 3.43368
$ws.Range("O20").Value = 0.007105781457704291
$ws.Range("P20").Value = 0.007218387680924443
$ws.Range("Q20").Value = 22.22402185976
$ws.Range("R20").Value = 200.01619673784
$ws.Range("S20").Value = 0.001824103199921536
$ws.Range("T20").Value = 0.00196153400991774

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 19.41708766666666
$ws.Range("H21").Value = 58.25126299999999
$ws.Range("I21").Value = 0.2567069098281641
$ws.Range("J21").Value = 0.2717412941260215
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 7.5382435
$ws.Range("N21").Value = 15.076487
$ws.Range("O21").Value = 0.04679974041200103
$ws.Range("P21").Value = 0.0316942545701456
$ws.Range("Q21").Value = 146.3707348921801
$ws.Range("R21").Value = 878.2244093530809
$ws.Range("S21").Value = 0.01201381674192504
$ws.Range("T21").Value = 0.008612637753250937

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 12.554453
$ws.Range("H22").Value = 25.108906
$ws.Range("I22").Value = 0.1659782810655758
$ws.Range("J22").Value = 0.1171326810635613
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 147.0592853333333
$ws.Range("N22").Value = 441.177856
$ws.Range("O22").Value = 0.9129893958419346
$ws.Range("P22").Value = 0.9274576550077637
$ws.Range("Q22").Value = 1846.248885930923
$ws.Range("R22").Value = 11077.49331558554
$ws.Range("S22").Value = 0.1515364105529429
$ws.Range("T22").Value = 0.1086356017039828

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 12.554453
$ws.Range("H23").Value = 25.108906
$ws.Range("I23").Value = 0.1659782810655758
$ws.Range("J23").Value = 0.1171326810635613
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.3688046666666667
$ws.Range("N23").Value = 1.106414
$ws.Range("O23").Value = 0.00228965310854373
$ws.Range("P23").Value = 0.002325937532793486
$ws.Range("Q23").Value = 4.630140853847333
$ws.Range("R23").Value = 27.780845123084
$ws.Range("S23").Value = 0.0003800326871925405
$ws.Range("T23").Value = 0.000272443299202466

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 12.554453
$ws.Range("H24").Value = 25.108906
$ws.Range("I24").Value = 0.1659782810655758
$ws.Range("J24").Value = 0.1171326810635613
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 4.963579
$ws.Range("N24").Value = 14.890737
$ws.Range("O24").Value = 0.03081542917981618
$ws.Range("P24").Value = 0.03130376520837289
$ws.Range("Q24").Value = 62.315019267287
$ws.Range("R24").Value = 373.8901156037219
$ws.Range("S24").Value = 0.005114691965563876
$ws.Range("T24").Value = 0.003666693946240946

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 12.554453
$ws.Range("H25").Value = 25.108906
$ws.Range("I25").Value = 0.1659782810655758
$ws.Range("J25").Value = 0.1171326810635613
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 1.14456
$ws.Range("N25").Value = 3.43368
$ws.Range("O25").Value = 0.007105781457704291
$ws.Range("P25").Value = 0.007218387680924443
$ws.Range("Q25").Value = 14.36932472568
$ws.Range("R25").Value = 86.21594835407998
$ws.Range("S25").Value = 0.0011794053919774
$ws.Range("T25").Value = 0.0008455091020228624

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 12.554453
$ws.Range("H26").Value = 25.108906
$ws.Range("I26").Value = 0.1659782810655758
$ws.Range("J26").Value = 0.1171326810635613
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 7.5382435
$ws.Range("N26").Value = 15.076487
$ws.Range("O26").Value = 0.04679974041200103
$ws.Range("P26").Value = 0.0316942545701456
$ws.Range("Q26").Value = 94.6385237233055
$ws.Range("R26").Value = 378.554094893222
$ws.Range("S26").Value = 0.007767740467899093
$ws.Range("T26").Value = 0.003712433012112183
